# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 257
    3  = 1332
    4  = 145
    7  = 94
    8  = 14
    9  = 179
    10 = 128
    11 = 4475
    12 = 6741
    18 = 4104
    19 = 477
    20 = 71
    21 = 51
    22 = 2687
    24 = 546
    26 = 350
    27 = 353
    28 = 395
    29 = 218
    31 = 1616
    32 = 1017
    33 = 60
    36 = 538
    40 = 630
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
